$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.364.17"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.27"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.77"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4475"
$ws.Range("E7").Value = "  -2.72%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3774"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07456"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8876"
$ws.Range("E10").Value = "  +2.93%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.99"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.97"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.746"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("E14").Value = "  +1.85%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.76"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07136"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008796"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.15"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.367.25"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.409"
$ws.Range("E22").Value = "  +3.87%  "

$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.972"
$ws.Range("E24").Value = "  -1.61%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.48"
$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.330"
$ws.Range("E26").Value = "  +4.64%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.67"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.398"
$ws.Range("E28").Value = "  +2.22%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.96"
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08906"
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7932"
$ws.Range("E31").Value = "  +3.28%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.603"
$ws.Range("E33").Value = "  +2.80%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.928"
$ws.Range("E34").Value = "  -1.51%  "

$ws.Range("B35").Value = "Frax"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.110"
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01990"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05308"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.304"
$ws.Range("E39").Value = "  +1.31%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5347"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.878"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1720"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.300"
$ws.Range("E43").Value = "  +15.86%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.668"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5108"
$ws.Range("E45").Value = "  -2.12%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  -0.88%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.697"
$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.24"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06412"
$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.04"
$ws.Range("E51").Value = "  +3.74%  "
